{"js": "// Replace the date line and every \"A\u00d7B=C\" multiplication answer in the\n// document with the values from the updated answer key. Each \"from\"\n// string is unique in the document, so a straightforward search +\n// replace per pair is unambiguous and order-independent.\nconst replacements = [\n  [\"2024-05-28 Tuesday\", \"2024-05-29 Wednesday\"],\n  [\"914\u00d72=1828\", \"775\u00d76=4650\"],\n  [\"138\u00d78=1104\", \"483\u00d72=966\"],\n  [\"994\u00d75=4970\", \"277\u00d77=1939\"],\n  [\"344\u00d77=2408\", \"288\u00d78=2304\"],\n  [\"253\u00d76=1518\", \"392\u00d77=2744\"],\n  [\"298\u00d78=2384\", \"580\u00d72=1160\"],\n  [\"359\u00d76=2154\", \"412\u00d75=2060\"],\n  [\"510\u00d72=1020\", \"725\u00d73=2175\"],\n  [\"132\u00d78=1056\", \"531\u00d75=2655\"],\n  [\"965\u00d75=4825\", \"425\u00d74=1700\"],\n  [\"264\u00d76=1584\", \"253\u00d75=1265\"],\n  [\"815\u00d75=4075\", \"638\u00d72=1276\"],\n  [\"401\u00d75=2005\", \"563\u00d73=1689\"],\n  [\"518\u00d77=3626\", \"608\u00d76=3648\"],\n  [\"746\u00d78=5968\", \"519\u00d76=3114\"],\n  [\"656\u00d79=5904\", \"523\u00d74=2092\"],\n  [\"779\u00d76=4674\", \"136\u00d73=408\"],\n  [\"947\u00d78=7576\", \"372\u00d72=744\"],\n  [\"576\u00d79=5184\", \"685\u00d77=4795\"],\n  [\"525\u00d79=4725\", \"489\u00d73=1467\"],\n  [\"895\u00d75=4475\", \"547\u00d77=3829\"],\n  [\"348\u00d74=1392\", \"115\u00d77=805\"],\n  [\"380\u00d78=3040\", \"596\u00d72=1192\"],\n  [\"570\u00d75=2850\", \"445\u00d79=4005\"],\n  [\"913\u00d79=8217\", \"211\u00d74=844\"],\n];\n\nconst body = context.document.body;\n\nfor (const [from, to] of replacements) {\n  const results = body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"A\u00d7B=C\" multiplication answer in the\n# document with the values from the updated answer key. Each \"from\"\n# string is unique in the document, so Find/Replace (first match,\n# MatchCase on) for each pair unambiguously targets the right run.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-05-28 Tuesday\", \"2024-05-29 Wednesday\"),\n    @(\"914\u00d72=1828\", \"775\u00d76=4650\"),\n    @(\"138\u00d78=1104\", \"483\u00d72=966\"),\n    @(\"994\u00d75=4970\", \"277\u00d77=1939\"),\n    @(\"344\u00d77=2408\", \"288\u00d78=2304\"),\n    @(\"253\u00d76=1518\", \"392\u00d77=2744\"),\n    @(\"298\u00d78=2384\", \"580\u00d72=1160\"),\n    @(\"359\u00d76=2154\", \"412\u00d75=2060\"),\n    @(\"510\u00d72=1020\", \"725\u00d73=2175\"),\n    @(\"132\u00d78=1056\", \"531\u00d75=2655\"),\n    @(\"965\u00d75=4825\", \"425\u00d74=1700\"),\n    @(\"264\u00d76=1584\", \"253\u00d75=1265\"),\n    @(\"815\u00d75=4075\", \"638\u00d72=1276\"),\n    @(\"401\u00d75=2005\", \"563\u00d73=1689\"),\n    @(\"518\u00d77=3626\", \"608\u00d76=3648\"),\n    @(\"746\u00d78=5968\", \"519\u00d76=3114\"),\n    @(\"656\u00d79=5904\", \"523\u00d74=2092\"),\n    @(\"779\u00d76=4674\", \"136\u00d73=408\"),\n    @(\"947\u00d78=7576\", \"372\u00d72=744\"),\n    @(\"576\u00d79=5184\", \"685\u00d77=4795\"),\n    @(\"525\u00d79=4725\", \"489\u00d73=1467\"),\n    @(\"895\u00d75=4475\", \"547\u00d77=3829\"),\n    @(\"348\u00d74=1392\", \"115\u00d77=805\"),\n    @(\"380\u00d78=3040\", \"596\u00d72=1192\"),\n    @(\"570\u00d75=2850\", \"445\u00d79=4005\"),\n    @(\"913\u00d79=8217\", \"211\u00d74=844\")\n)\n\nforeach ($pair in $pairs) {\n    $from = $pair[0]\n    $to = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($from, $false, $true, $false, $false, $false, $true, 1, $false, $to, 2)\n}\n"}
